# "Got all BOMs up to date" — refresh the BOM sheet:
#  - update quantities/unit-prices/subtotals on several existing rows
#  - insert two new component rows (B4B-XH-A connector, DX07S024WJ3R400 connector)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update existing rows whose Quantity / Unit Price / Subtotal changed
# ---------------------------------------------------------------------------

# Row 4 - RC0603JR-07470RL (470 Res1): unit price / subtotal
$ws.Range("G4").Value = 0.094
$ws.Range("H4").Value = 0.094

# Row 5 - CR0603-JW-331ELF (330 Res1): unit price / subtotal
$ws.Range("G5").Value = 0.0028
$ws.Range("H5").Value = 0.07

# Row 6 - CL10B104KA8NNNC (100nF Cap): quantity + subtotal
$ws.Range("D6").Value = 56
$ws.Range("H6").Value = 0.8176

# Row 9 - RMCF0603FT13K7 (13.7k Res1): quantity + subtotal
$ws.Range("D9").Value = 13
$ws.Range("H9").Value = 0.195

# Row 13 - 0603B472K500CT (4.7nF Cap): quantity, unit price, subtotal
$ws.Range("D13").Value = 7
$ws.Range("G13").Value = 0.02163
$ws.Range("H13").Value = 0.21628

# Row 14 - RC0603JR07-4K7L (4.7k Res1): unit price + subtotal
$ws.Range("G14").Value = 0.015
$ws.Range("H14").Value = 0.27

# Row 15 - GRM1885C1H102JA01D (1nF Cap): unit price + subtotal
$ws.Range("G15").Value = 0.0483
$ws.Range("H15").Value = 0.1449

# Row 16 - RC0603JR-070RL (0 Res1): unit price + subtotal
$ws.Range("G16").Value = 0.0158
$ws.Range("H16").Value = 0.0158

# ---------------------------------------------------------------------------
# 2) Insert the new "B4B-XH-A (LF)(SN)" connector row as new row 17
#    (pushes D5V0L2B3W-7 and everything after it down by one)
# ---------------------------------------------------------------------------
$ws.Range("A17").EntireRow.Insert()

$ws.Range("A17").Value = "B4B-XH-A (LF)(SN)"
$ws.Range("B17").Value = "No Description Available"
$ws.Range("C17").Value = "B4B-XH-A(LF)(SN)"
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = "CONN_B4B-XH-A (LF)(SN)_JST"
$ws.Range("G17").Value = 0.21
$ws.Range("H17").Value = 0.21

# ---------------------------------------------------------------------------
# 3) Insert the new "DX07S024WJ3R400" connector row as new row 19
#    (after D5V0L2B3W-7, which is now row 18; pushes Header 5 etc. down by one)
# ---------------------------------------------------------------------------
$ws.Range("A19").EntireRow.Insert()

$ws.Range("A19").Value = "DX07S024WJ3R400"
$ws.Range("B19").Value = "No Description Available"
$ws.Range("C19").Value = "DX07S024WJ3R400"
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = "CONN_DX07S024WJ1R350_JAE"
$ws.Range("G19").Value = 2.11
$ws.Range("H19").Value = 2.11

# ---------------------------------------------------------------------------
# 4) Update rows that shifted down by two (old 20/21/22 -> new 22/23/24)
# ---------------------------------------------------------------------------

# Row 22 (was LSM6DSMTR @ old row 20): unit price + subtotal
$ws.Range("G22").Value = 2.9
$ws.Range("H22").Value = 2.9

# Row 23 (was MCP2515-I/SO @ old row 21): unit price + subtotal
$ws.Range("G23").Value = 1.71
$ws.Range("H23").Value = 1.71

# Row 24 (was MP5087GG @ old row 22): quantity + subtotal
$ws.Range("D24").Value = 13
$ws.Range("H24").Value = 14.17
